# The document repeats the same "campaign period" sentence in 4 places.
# Each instance is built from a long chain of differently-formatted runs
# (per-character / per-word Meiryo runs with w:sz=18, plus a couple of
# w:proofErr gramStart/gramEnd markers around "ペルセウス"). The edit
# collapses each of those run chains into a single, plain run (no rPr at
# all) containing the new translated/updated sentence about Cygnus.

$d = $word.ActiveDocument

$oldMarker = "年キャンペーン期間"
$newText = "年キャンペーン期間 対象：Cygnus: 8月10〜19日、9月9〜18日、10月8〜18日"

# Collect the paragraphs that hold the old sentence first (by index),
# since deleting/inserting content shifts character offsets but not
# paragraph identity/order.
$count = $d.Paragraphs.Count
$targets = @()
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$oldMarker*") {
        $targets += $i
    }
}

# Walk backwards so earlier paragraph indices stay valid while we edit.
for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $p = $d.Paragraphs($targets[$j])

    # Range over the paragraph's content, excluding the trailing
    # paragraph mark, covering every run (plain text + proofErr markers).
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null

    # Delete all existing (variously formatted) runs, then insert fresh,
    # unformatted text -- this yields a single <w:r><w:t>...</w:t></w:r>
    # with no <w:rPr>, matching the target edit exactly.
    $r.Delete()

    $r2 = $p.Range
    $r2.MoveEnd(1, -1) | Out-Null
    $r2.InsertAfter($newText)
}

Write-Host "Replaced" $targets.Count "occurrence(s) of the campaign-period sentence."
